$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.615.91'
$ws.Range('E2').Value = '  -0.45%  '

$ws.Range('D3').Value = '3.472.02'
$ws.Range('E3').Value = '  +0.52%  '

$ws.Range('D5').Value = '575.17'
$ws.Range('E5').Value = '  -0.17%  '

$ws.Range('D6').Value = '160.68'
$ws.Range('E6').Value = '  +0.13%  '

$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.16%  '

$ws.Range('D8').Value = '3.474.63'
$ws.Range('E8').Value = '  +0.54%  '

$ws.Range('D9').Value = '0.574'
$ws.Range('E9').Value = '  -5.71%  '

$ws.Range('E10').Value = '  -0.36%  '

$ws.Range('E11').Value = '  -2.31%  '

$ws.Range('E12').Value = '  -3.59%  '

$ws.Range('D13').Value = '4.062.07'
$ws.Range('E13').Value = '  +0.36%  '

$ws.Range('E14').Value = '  -0.17%  '

$ws.Range('D15').Value = '27.61'
$ws.Range('E15').Value = '  -2.28%  '

$ws.Range('E16').Value = '  -8.39%  '

$ws.Range('D17').Value = '64.659.42'
$ws.Range('E17').Value = '  -0.41%  '

$ws.Range('D18').Value = '3.528.88'
$ws.Range('E18').Value = '  +0.97%  '

$ws.Range('D19').Value = '6.25'
$ws.Range('E19').Value = '  -3.44%  '

$ws.Range('E20').Value = '  -3.26%  '

$ws.Range('D21').Value = '381.95'
$ws.Range('E21').Value = '  +0.26%  '

$ws.Range('D22').Value = '7.96'
$ws.Range('E22').Value = '  -2.51%  '

$ws.Range('D23').Value = '72.60'
$ws.Range('E23').Value = '  -0.51%  '

$ws.Range('E24').Value = '  +0.25%  '

$ws.Range('D25').Value = '0.529'
$ws.Range('E25').Value = '  -4.61%  '

$ws.Range('E26').Value = '  -0.67%  '

$ws.Range('E27').Value = '  -1.82%  '

$ws.Range('E28').Value = '  +0.88%  '

$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.10%  '

$ws.Range('E30').Value = '  -0.30%  '

$ws.Range('E31').Value = '  -6.06%  '

$ws.Range('E32').Value = '  -1.46%  '

$ws.Range('D33').Value = '23.35'
$ws.Range('E33').Value = '  -1.06%  '

$ws.Range('E34').Value = '  -2.88%  '

$ws.Range('D35').Value = '1.59'
$ws.Range('E35').Value = '  -1.30%  '

$ws.Range('D36').Value = '160.85'
$ws.Range('E36').Value = '  -0.15%  '

$ws.Range('E37').Value = '  -2.89%  '

$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').Value = '26.98'
$ws.Range('E38').Value = '  +1.64%  '

$ws.Range('B39').Value = 'Mantle'
$ws.Range('C39').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D39').Value = '0.823'
$ws.Range('E39').Value = '  +5.57%  '

$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '0.0748'
$ws.Range('E40').Value = '  -3.81%  '

$ws.Range('D41').Value = '2.844.44'
$ws.Range('E41').Value = '  -2.36%  '

$ws.Range('D42').Value = '4.51'
$ws.Range('E42').Value = '  -3.11%  '

$ws.Range('D43').Value = '42.89'
$ws.Range('E43').Value = '  -0.28%  '

$ws.Range('D44').Value = '6.46'
$ws.Range('E44').Value = '  -5.09%  '

$ws.Range('D45').Value = '25.85'
$ws.Range('E45').Value = '  -0.55%  '

$ws.Range('D46').Value = '0.0309'
$ws.Range('E46').Value = '  -3.22%  '

$ws.Range('E47').Value = '  +10.79%  '

$ws.Range('D48').Value = '336.04'
$ws.Range('E48').Value = '  +4.17%  '

$ws.Range('E49').Value = '  -2.61%  '

$ws.Range('E50').Value = '  -2.10%  '

$ws.Range('E51').Value = '  -3.47%  '
